$d = $word.ActiveDocument

$oldText = "Connor: Has previous experience in computer stuff and some programming languages?"

# Locate the paragraph that still contains the original sentence (literal
# substring match -- avoid -like since "?" is a wildcard there).
$target = $null
for ($i = 1; $i -le $d.Paragraphs.Count; $i++) {
    $p = $d.Paragraphs($i)
    if ($p.Range.Text.Contains($oldText)) {
        $target = $p
        break
    }
}

$start = $target.Range.Start
$end = $start + $oldText.Length

# Remove the old single run's text, then insert the replacement as a set of
# distinct <w:r> elements (one per diff hunk) via InsertXML so they are not
# coalesced into a single run on save.
$oldRange = $d.Range($start, $end)
$oldRange.Text = ""

$insertPoint = $d.Range($start, $start)

$rightQuote = [char]0x2019

$xml = '<?xml version="1.0" encoding="UTF-8" standalone="yes"?>' +
'<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage">' +
'<pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml">' +
'<pkg:xmlData>' +
'<w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main">' +
'<w:body>' +
'<w:p>' +
'<w:r><w:t xml:space="preserve">Connor: Has previous experience in </w:t></w:r>' +
('<w:r><w:t xml:space="preserve">building PC' + $rightQuote + 's </w:t></w:r>') +
'<w:r><w:t>and general hardware troubleshooting, as well as dabbling in some cod</w:t></w:r>' +
'<w:r><w:t>ing languages</w:t></w:r>' +
'<w:r><w:t xml:space="preserve"> such as python</w:t></w:r>' +
'<w:r><w:t>.</w:t></w:r>' +
'<w:r><w:t xml:space="preserve"> </w:t></w:r>' +
'</w:p>' +
'</w:body>' +
'</w:document>' +
'</pkg:xmlData>' +
'</pkg:part>' +
'</pkg:package>'

$insertPoint.InsertXML($xml)
